$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rows = @(
    ,@("MuSCs", "Hc", "C5ar2", "ECs", 1, 0.3333333333333333, 0.015297, 0.045891, 0.0153746470927005, 0.0153746470927005, 2, 0.6666666666666666, 0.007121666666666668, 0.021365, 0.00006626332888101038, 0.00006626332888101038, 0.000108940135, 0.0009804612150000002, 0.000001018775296733084, 0.000001018775296733084)
    ,@("MuSCs", "Hc", "C5ar2", "FAPs", 1, 0.3333333333333333, 0.015297, 0.045891, 0.0153746470927005, 0.0153746470927005, 1, 0.3333333333333333, 0.01171, 0.03513, 0.0001089553355295996, 0.0001089553355295996, 0.00017912787, 0.00161215083, 0.000001675149832634366, 0.000001675149832634366)
    ,@("MuSCs", "Hc", "C5ar2", "Inflammatory-Mac", 1, 0.3333333333333333, 0.015297, 0.045891, 0.0153746470927005, 0.0153746470927005, 3, 1, 15.223647, 45.67094100000001, 0.141647956180118, 0.141647956180118, 0.232876128159, 2.095885153431, 0.002177787337671619, 0.002177787337671619)
    ,@("MuSCs", "Hc", "C5ar2", "MuSCs", 1, 0.3333333333333333, 0.015297, 0.045891, 0.0153746470927005, 0.0153746470927005, 1, 0.3333333333333333, 0.01159, 0.03477, 0.000107838799213327, 0.000107838799213327, 0.00017729223, 0.00159563007, 0.000001657983480805491, 0.000001657983480805491)
    ,@("MuSCs", "Hc", "C5ar2", "Neutrophils", 1, 0.3333333333333333, 0.015297, 0.045891, 0.0153746470927005, 0.0153746470927005, 3, 1, 84.97145566666667, 254.914367, 0.7906143008198258, 0.7906143008198258, 1.299808357333, 11.698275215997, 0.01215541586154698, 0.01215541586154698)
    ,@("MuSCs", "Hc", "C5ar2", "Resolving-Mac", 1, 0.3333333333333333, 0.015297, 0.045891, 0.0153746470927005, 0.0153746470927005, 3, 1, 7.249708000000001, 21.749124, 0.0674546855364323, 0.06745468553643229, 0.110898783276, 0.9980890494840001, 0.001037091984871736, 0.001037091984871735)
    ,@("Neutrophils", "Hc", "C5ar2", "ECs", 3, 1, 0.9796526666666665, 2.938958, 0.9846253529072995, 0.9846253529072995, 2, 0.6666666666666666, 0.007121666666666668, 0.021365, 0.00006626332888101038, 0.00006626332888101038, 0.006976759741111111, 0.06279083767, 0.0000652445535842773, 0.0000652445535842773)
    ,@("Neutrophils", "Hc", "C5ar2", "FAPs", 3, 1, 0.9796526666666665, 2.938958, 0.9846253529072995, 0.9846253529072995, 1, 0.3333333333333333, 0.01171, 0.03513, 0.0001089553355295996, 0.0001089553355295996, 0.01147173272666666, 0.10324559454, 0.0001072801856969652, 0.0001072801856969652)
    ,@("Neutrophils", "Hc", "C5ar2", "Inflammatory-Mac", 3, 1, 0.9796526666666665, 2.938958, 0.9846253529072995, 0.9846253529072995, 3, 1, 15.223647, 45.67094100000001, 0.141647956180118, 0.141647956180118, 14.913886379942, 134.224977419478, 0.1394701688424464, 0.1394701688424464)
    ,@("Neutrophils", "Hc", "C5ar2", "MuSCs", 3, 1, 0.9796526666666665, 2.938958, 0.9846253529072995, 0.9846253529072995, 1, 0.3333333333333333, 0.01159, 0.03477, 0.000107838799213327, 0.000107838799213327, 0.01135417440666667, 0.10218756966, 0.0001061808157325215, 0.0001061808157325215)
    ,@("Neutrophils", "Hc", "C5ar2", "Neutrophils", 3, 1, 0.9796526666666665, 2.938958, 0.9846253529072995, 0.9846253529072995, 3, 1, 84.97145566666667, 254.914367, 0.7906143008198258, 0.7906143008198258, 83.24251313439844, 749.1826182095858, 0.7784588849582788, 0.7784588849582788)
    ,@("Neutrophils", "Hc", "C5ar2", "Resolving-Mac", 3, 1, 0.9796526666666665, 2.938958, 0.9846253529072995, 0.9846253529072995, 3, 1, 7.249708000000001, 21.749124, 0.0674546855364323, 0.06745468553643229, 7.102195774754667, 63.919761972792, 0.06641759355156056, 0.06641759355156054)
)

$startRow = 2
for ($i = 0; $i -lt $rows.Count; $i++) {
    $r = $startRow + $i
    $rowData = $rows[$i]
    for ($c = 1; $c -le $rowData.Count; $c++) {
        $ws.Cells.Item($r, $c).Value = $rowData[$c - 1]
    }
}
